$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and 1h volume change (E) values for each coin row.
# Values that look like plain numbers must be forced to remain text (to preserve
# formatting such as trailing zeros / multi-dot separators), matching the source data.

$ws.Range('D2').Value = '61.247.85'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '2.394.55'
$ws.Range('E3').Value = '  -3.70%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'549.83"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = "'142.09"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -10.57%  '
$ws.Range('D9').Value = '2.394.88'
$ws.Range('E9').Value = '  -3.69%  '
$ws.Range('E10').Value = '  -2.10%  '
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = "'5.30"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.83%  '
$ws.Range('E13').Value = '  -3.04%  '
$ws.Range('D14').Value = "'25.55"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.81%  '
$ws.Range('D15').Value = '2.826.49'
$ws.Range('E15').Value = '  -3.78%  '
$ws.Range('D16').Value = "'0.0000166"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').Value = '60.678.44'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').Value = '2.393.72'
$ws.Range('E18').Value = '  -3.95%  '
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').Value = "'319.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('E22').Value = '  -4.43%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = "'1.93"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('D25').Value = "'63.70"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').Value = "'8.27"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.91%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '2.512.84'
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('D29').Value = '0.0₃0937'
$ws.Range('E29').Value = '  -6.49%  '
$ws.Range('D30').Value = "'531.62"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.21%  '
$ws.Range('D31').Value = "'1.44"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.66%  '
$ws.Range('D32').Value = "'8.12"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.90%  '
$ws.Range('D33').Value = "'0.146"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.64%  '
$ws.Range('D34').Value = "'1.86"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.93%  '
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -2.93%  '
$ws.Range('D38').Value = "'5.57"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.79%  '
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('E40').Value = '  +6.58%  '
$ws.Range('D41').Value = "'18.16"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.39%  '
$ws.Range('D42').Value = "'138.60"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.76%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = "'40.29"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('E45').Value = '  -8.76%  '
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = "'141.33"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.76%  '
$ws.Range('E48').Value = '  -4.63%  '
$ws.Range('D49').Value = "'0.0521"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.63%  '
$ws.Range('D50').Value = "'0.579"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.44%  '
$ws.Range('D51').Value = "'0.0227"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.09%  '
